$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column string values (prices with "." as thousands separators)
# are written as text, not auto-converted to numbers, then restore the
# original (unstyled) cell style so no formatting diff is introduced.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "53.972.10"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "2.279.81"
$ws.Range("E3").Value = "  +1.97%  "
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "494.38"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").Value = "127.21"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").Value = "0.994"
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("E8").Value = "  +1.67%  "
$ws.Range("D9").Value = "2.277.42"
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("D10").Value = "0.0948"
$ws.Range("E10").Value = "  +3.19%  "
$ws.Range("E11").Value = "  +2.15%  "
$ws.Range("D12").Value = "0.325"
$ws.Range("E12").Value = "  +3.27%  "
$ws.Range("E13").Value = "  -2.32%  "
$ws.Range("D14").Value = "2.652.82"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").Value = "21.79"
$ws.Range("E15").Value = "  +3.19%  "
$ws.Range("D16").Value = "53.926.69"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("D18").Value = "2.252.43"
$ws.Range("E18").Value = "  +4.67%  "
$ws.Range("E19").Value = "  +4.70%  "
$ws.Range("D20").Value = "4.10"
$ws.Range("E20").Value = "  +3.56%  "
$ws.Range("D21").Value = "6.45"
$ws.Range("E21").Value = "  +5.38%  "
$ws.Range("D22").Value = "300.85"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("D23").Value = "0.992"
$ws.Range("E23").Value = "  -0.90%  "
$ws.Range("D24").Value = "5.38"
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("D25").Value = "62.51"
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +0.72%  "
$ws.Range("E27").Value = "  +2.29%  "
$ws.Range("E28").Value = "  +4.23%  "
$ws.Range("D29").Value = "2.347.00"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").Value = "7.07"
$ws.Range("E30").Value = "  +0.74%  "
$ws.Range("D31").Value = "168.70"
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("D32").Value = "1.60"
$ws.Range("E32").Value = "  +0.77%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("E34").Value = "  +1.88%  "
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").Value = "0.996"
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("D37").Value = "1.07"
$ws.Range("E37").Value = "  +1.64%  "
$ws.Range("D38").Value = "17.58"
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("D39").Value = "1.19"
$ws.Range("E39").Value = "  +2.89%  "
$ws.Range("D40").Value = "0.862"
$ws.Range("E40").Value = "  +2.98%  "
$ws.Range("D41").Value = "3.70"
$ws.Range("E41").Value = "  +3.81%  "
$ws.Range("D42").Value = "35.36"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("E43").Value = "  +2.96%  "
$ws.Range("E44").Value = "  +2.31%  "
$ws.Range("D45").Value = "3.35"
$ws.Range("E45").Value = "  +1.55%  "
$ws.Range("D46").Value = "127.87"
$ws.Range("E46").Value = "  +4.05%  "
$ws.Range("D47").Value = "4.91"
$ws.Range("E47").Value = "  +5.86%  "
$ws.Range("D48").Value = "0.0888"
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("D50").Value = "238.16"
$ws.Range("E50").Value = "  +2.81%  "
$ws.Range("D51").Value = "0.0483"
$ws.Range("E51").Value = "  +2.77%  "

$dRange.Style = "Normal"

